$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits to existing cells (label tweaks) ---
$ws.Range("A3").Value = "Health indicator, Condition indicator, Could be a difference between the two"
$ws.Range("A5").Value = "Processing of measured variable"
$ws.Range("A7").Value = "Heath state- Health indicator mapping, fk"
$ws.Range("A8").Value = "Health indicator - Measured variable mapping, hk"

# --- New column E content (order-spectrum / TSA order-spectrum notes) ---
$ws.Range("E2").Value = "Spiral bevel gear degradation (pitting)"
$ws.Range("E3").Value = "Oil debris mass"
$ws.Range("E4").Value = "Acceleration, Oli debris mass"
$ws.Range("E5").Value = "One dimensional transition function using whitening transform, TSA many other CI's."
$ws.Range("E7").Value = "Direct"
$ws.Range("E8").Value = " Data driven Double exponential smoothing model"
$ws.Range("E10").Value = "N/A This was buildt on data (ARIMA)"
$ws.Range("E11").Value = "Particle Filter with l-step ahead estimator"

# --- New row 14: flagged "Bad" style note cell ---
$ws.Range("E14").Value = "Check if this is summarized in lit review"
$ws.Range("E14").Style = "Bad"

# --- Row 3 grew taller to fit the new wrapped text ---
$ws.Range("A3:N3").RowHeight = 60

# --- Selection / scroll position moved to the newly edited cell ---
$ws.Range("E4").Select()
